# Update numeric profit figures across multiple sheets as part of the scheduled
# Gungnir_Profits refresh. Values are written with explicit .Value assignments so
# that exact numeric literals (including floating point values) are preserved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 13889421
$ws.Range("I2").Value = 20833382
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 20833382
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -20833269
$ws.Range("N2").Value = -1726
$ws.Range("H9").Value = 67.90000000000001
$ws.Range("I9").Value = 82.375
$ws.Range("J9").Value = 10
$ws.Range("K9").Value = 82.375
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = 86.625
$ws.Range("N9").Value = -348
$ws.Range("H43").Value = 1543.2285
$ws.Range("I43").Value = 1725.7858
$ws.Range("J43").Value = 1421.5238
$ws.Range("K43").Value = 1725.7858
$ws.Range("L43").Value = 1421.5238
$ws.Range("M43").Value = -1656.7858
$ws.Range("N43").Value = -1559.5238
$ws.Range("H86").Value = 6120606
$ws.Range("I86").Value = 100000
$ws.Range("K86").Value = 100000
$ws.Range("M86").Value = -98877
$ws.Range("H89").Value = 6120606
$ws.Range("I89").Value = 100000
$ws.Range("K89").Value = 500000
$ws.Range("M89").Value = -494384
$ws.Range("H112").Value = 12135.889
$ws.Range("J112").Value = 12466.914
$ws.Range("L112").Value = 37400.742
$ws.Range("N112").Value = -39616.742
$ws.Range("H135").Value = 1053.5588
$ws.Range("I135").Value = 1024.8788
$ws.Range("K135").Value = 9223.9092
$ws.Range("M135").Value = -6688.9092
$ws.Range("H138").Value = 1738.4828
$ws.Range("I138").Value = 1055.9773
$ws.Range("J138").Value = 3883.5
$ws.Range("K138").Value = 3167.9319
$ws.Range("L138").Value = 11650.5
$ws.Range("M138").Value = 1972.0681
$ws.Range("N138").Value = -21930.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1127.73
$ws.Range("I32").Value = 1044.979
$ws.Range("J32").Value = 2700
$ws.Range("K32").Value = 1044.979
$ws.Range("L32").Value = 2700
$ws.Range("M32").Value = -757.979
$ws.Range("N32").Value = -3274
$ws.Range("H74").Value = 1410.2894
$ws.Range("I74").Value = 1460.3636
$ws.Range("J74").Value = 1079.8
$ws.Range("K74").Value = 1460.3636
$ws.Range("L74").Value = 1079.8
$ws.Range("M74").Value = -586.3635999999999
$ws.Range("N74").Value = -2827.8
$ws.Range("H77").Value = 1410.2894
$ws.Range("I77").Value = 1460.3636
$ws.Range("J77").Value = 1079.8
$ws.Range("K77").Value = 7301.817999999999
$ws.Range("L77").Value = 5399
$ws.Range("M77").Value = -2933.817999999999
$ws.Range("N77").Value = -14135
$ws.Range("H122").Value = 729.875
$ws.Range("I122").Value = 736.5
$ws.Range("J122").Value = 710
$ws.Range("K122").Value = 2209.5
$ws.Range("L122").Value = 2130
$ws.Range("M122").Value = 240.5
$ws.Range("N122").Value = -7030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20834298
$ws.Range("I58").Value = 32258800
$ws.Range("J58").Value = 1384
$ws.Range("K58").Value = 32258800
$ws.Range("L58").Value = 1384
$ws.Range("M58").Value = -32258597
$ws.Range("N58").Value = -1790
$ws.Range("H68").Value = 21397.5
$ws.Range("I68").Value = 20000
$ws.Range("K68").Value = 20000
$ws.Range("M68").Value = -19251
$ws.Range("H71").Value = 21397.5
$ws.Range("I71").Value = 20000
$ws.Range("K71").Value = 60000
$ws.Range("M71").Value = -56256
$ws.Range("H74").Value = 41264
$ws.Range("J74").Value = 41264
$ws.Range("L74").Value = 41264
$ws.Range("N74").Value = -43012
$ws.Range("H77").Value = 41264
$ws.Range("J77").Value = 41264
$ws.Range("L77").Value = 123792
$ws.Range("N77").Value = -132528
$ws.Range("H132").Value = 11495868
$ws.Range("I132").Value = 1329.0555
$ws.Range("J132").Value = 30305112
$ws.Range("K132").Value = 3987.1665
$ws.Range("L132").Value = 90915336
$ws.Range("M132").Value = -1457.1665
$ws.Range("N132").Value = -90920396
$ws.Range("H136").Value = 20834298
$ws.Range("I136").Value = 32258800
$ws.Range("J136").Value = 1384
$ws.Range("K136").Value = 96776400
$ws.Range("L136").Value = 4152
$ws.Range("M136").Value = -96773850
$ws.Range("N136").Value = -9252

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 37638284
$ws.Range("I5").Value = 51282536
$ws.Range("J5").Value = 27784100
$ws.Range("K5").Value = 153847608
$ws.Range("L5").Value = 83352300
$ws.Range("M5").Value = -153847496
$ws.Range("N5").Value = -83352524
$ws.Range("H113").Value = 17882384
$ws.Range("I113").Value = 9259663
$ws.Range("J113").Value = 21256492
$ws.Range("K113").Value = 27778989
$ws.Range("L113").Value = 63769476
$ws.Range("M113").Value = -27776819
$ws.Range("N113").Value = -63773816
$ws.Range("H122").Value = 9770210
$ws.Range("I122").Value = 56818452
$ws.Range("J122").Value = 5479.7925
$ws.Range("K122").Value = 511366068
$ws.Range("L122").Value = 49318.13249999999
$ws.Range("M122").Value = -511363618
$ws.Range("N122").Value = -54218.13249999999
$ws.Range("H131").Value = 14961297
$ws.Range("I131").Value = 51283220
$ws.Range("J131").Value = 7696911.5
$ws.Range("K131").Value = 153849660
$ws.Range("L131").Value = 23090734.5
$ws.Range("M131").Value = -153844620
$ws.Range("N131").Value = -23100814.5
$ws.Range("H135").Value = 37638284
$ws.Range("I135").Value = 51282536
$ws.Range("J135").Value = 27784100
$ws.Range("K135").Value = 461542824
$ws.Range("L135").Value = 250056900
$ws.Range("M135").Value = -461540289
$ws.Range("N135").Value = -250061970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6253294.5
$ws.Range("I80").Value = 3550.7144
$ws.Range("K80").Value = 3550.7144
$ws.Range("M80").Value = -2552.7144
$ws.Range("H83").Value = 6253294.5
$ws.Range("I83").Value = 3550.7144
$ws.Range("K83").Value = 17753.572
$ws.Range("M83").Value = -12761.572
$ws.Range("H132").Value = 8272.714
$ws.Range("I132").Value = 5418.2856
$ws.Range("J132").Value = 16836
$ws.Range("K132").Value = 16254.8568
$ws.Range("L132").Value = 50508
$ws.Range("M132").Value = -13724.8568
$ws.Range("N132").Value = -55568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20414198
$ws.Range("I132").Value = 33334972
$ws.Range("J132").Value = 12973.842
$ws.Range("K132").Value = 100004916
$ws.Range("L132").Value = 38921.526
$ws.Range("M132").Value = -100002386
$ws.Range("N132").Value = -43981.526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 30389
$ws.Range("J93").Value = 30389
$ws.Range("L93").Value = 30389
$ws.Range("N93").Value = -35381
$ws.Range("H132").Value = 6039160
$ws.Range("I132").Value = 17649.5
$ws.Range("K132").Value = 52948.5
$ws.Range("M132").Value = -50418.5
$ws.Range("H136").Value = 9620004
$ws.Range("I136").Value = 11909659
$ws.Range("J136").Value = 3456.9
$ws.Range("K136").Value = 35728977
$ws.Range("L136").Value = 10370.7
$ws.Range("M136").Value = -35726427
$ws.Range("N136").Value = -15470.7
